# Rename the worksheet from "Domino Qi Mini Rev. D" to
# "Domino Qi Mini Rev. E" (board revision bump to a 4-layer design).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Domino Qi Mini Rev. D")
$oldName = $ws.Name
$newName = "Domino Qi Mini Rev. E"
$ws.Name = $newName

# Renaming the sheet automatically updates most of the workbook-scoped
# defined names (the numerous Print_Area_0, Print_Area_0_0, ... aliases)
# because their formulas are re-evaluated against the sheet object. The
# single plain "_xlnm.Print_Area" name however keeps its old literal
# sheet-name text in its RefersTo formula, so patch it explicitly too.
foreach ($n in $wb.Names) {
    if ($n.RefersTo -like "*'$oldName'*") {
        $n.RefersTo = $n.RefersTo -replace [regex]::Escape("'$oldName'"), "'$newName'"
    }
}
